$wb = $excel.ActiveWorkbook

# Update the "想去人数" (want-to-go count) column F on both the "展览"
# sheet and the duplicate "全部类型" sheet, which mirror the same data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 332
    $ws.Range("F3").Value = 240
    $ws.Range("F4").Value = 67
}
